# Q4 2022 Fiscal Data update
# Adds a new RMO row (row 37) to Sheet1 with the latest BIR collection
# goal figures (RMO No. 53-2022, issued December 7, 2022).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Write the new row's text values first in the same order Excel would
# append them to the shared-strings table (RMO, IssueDate, then Link),
# followed by the TargetCY / Version values that already exist in the
# shared-strings table, then the numeric goal figures.
$ws.Range("B37").Value = "RMO No. 53-2022"
$ws.Range("C37").Value = "December 7, 2022"
$ws.Range("D37").Value = "2022"
$ws.Range("E37").Value = "Final"
$ws.Range("A37").Value = "https://www.bir.gov.ph/images/bir_files/internal_communications_3/2022/Full%20Text/Revised%20CY2022%20RMO%20Goal%20Annexes.pdf"

$ws.Range("F37").Value = 2392587
$ws.Range("G37").Value = 1197966
$ws.Range("H37").Value = 365197
$ws.Range("I37").Value = 430160.553
$ws.Range("J37").Value = 153695.364
$ws.Range("K37").Value = 245568

# Update the frozen-pane scroll position / active cell to reflect the
# newly added row, matching the author's saved view state.
$ws.Range("A38").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
